# Table for I-V.xlsx edit script
# Rebuilds the I-V data table: relabels a header, clears/repopulates the
# "light power"/I_SC/V_OC measurement rows with real experimental values,
# and re-derives the Fill-factor / Conversion-efficiency formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Header row -----------------------------------------------------
$ws.Range("B1").Value2 = "Light power"
$ws.Range("D1").Value2 = "I_SC"
$ws.Range("A1").EntireRow.RowHeight = 29

# ---- Row 2 (Ref) : drop the stale Opti. Power formula ----------------
$ws.Range("E2").ClearContents()

# ---- Row 3 (SC1-Z907) : real measured values --------------------------
$ws.Range("B3").Value2 = 1.03
$ws.Range("C3").Value2 = [double]"0.43555094523216797"
$ws.Range("D3").Value2 = [double]"5.95059E-4"
$ws.Range("E3").Value2 = [double]"1.4060661561480001E-4"
$ws.Range("F3").Formula = "=E3/B3"
$ws.Range("G3").Formula = "=E3/(C3*D3)"

# logger-pasted cells: small monospace font, black text, left/centre
foreach ($addr in @("C3", "E3")) {
    $rr = $ws.Range($addr)
    $rr.Font.Name = "Courier New"
    $rr.Font.Size = 7
    $rr.Font.Color = 0
    $rr.HorizontalAlignment = -4131
    $rr.VerticalAlignment = -4108
}

# ---- Row 4 (SC2-Z907) : no data collected, keep error formulas --------
$ws.Range("E4").ClearContents()
$ws.Range("F4").Formula = "=E4/B4"
$ws.Range("G4").Formula = "=E4/(C4*D4)"

# ---- Row 5 (SC3-Tur) --------------------------------------------------
$ws.Range("E5").ClearContents()
$ws.Range("B5").Value2 = 1.03
$ws.Range("C5").Value2 = [double]"-6.91728065845537E-2"
$ws.Range("D5").Value2 = [double]"1.8923982E-7"
$ws.Range("F5").Formula = "=E5/B5"
$ws.Range("G5").Formula = "=E5/(C5*D5)"

# out-of-range reading: whole row flagged red
$ws.Range("B5").Font.Color = 255
$rr = $ws.Range("C5")
$rr.Font.Name = "Courier New"
$rr.Font.Size = 7
$rr.Font.Color = 255
$rr.HorizontalAlignment = -4131
$rr.VerticalAlignment = -4108
$rr = $ws.Range("D5")
$rr.NumberFormat = "0.00E+00"
$rr.Font.Color = 255

# ---- Row 6 (SC4-BeetR) -------------------------------------------------
$ws.Range("E6").ClearContents()
$ws.Range("B6").Value2 = 1.03
$ws.Range("C6").Value2 = [double]"-5.0010018548270398E-2"
$ws.Range("D6").Value2 = [double]"3.1655643999999998E-6"
$ws.Range("F6").Formula = "=E6/B6"
$ws.Range("G6").Formula = "=E6/(C6*D6)"

$ws.Range("B6").Font.Color = 255
$rr = $ws.Range("C6")
$rr.Font.Name = "Courier New"
$rr.Font.Size = 7
$rr.Font.Color = 255
$rr.HorizontalAlignment = -4131
$rr.VerticalAlignment = -4108
$rr = $ws.Range("D6")
$rr.NumberFormat = "0.00E+00"
$rr.Font.Color = 255

# ---- Row 7 (SC5-BlueB) -------------------------------------------------
$ws.Range("E7").ClearContents()
$ws.Range("B7").Value2 = 1.03
$ws.Range("C7").Value2 = [double]"0.13678308581395601"
$ws.Range("D7").Value2 = [double]"1.118747E-5"
$ws.Range("E7").Value2 = [double]"5.5025570662209095E-7"
$ws.Range("F7").Formula = "=E7/B7"
$ws.Range("G7").Formula = "=E7/(C7*D7)"

$rr = $ws.Range("C7")
$rr.Font.Name = "Courier New"
$rr.Font.Size = 7
$rr.Font.Color = 0
$rr.HorizontalAlignment = -4131
$rr.VerticalAlignment = -4108
$ws.Range("D7").NumberFormat = "0.00E+00"
$rr = $ws.Range("E7")
$rr.NumberFormat = "0.00E+00"
$rr.Font.Name = "Courier New"
$rr.Font.Size = 7
$rr.Font.Color = 0
$rr.HorizontalAlignment = -4131
$rr.VerticalAlignment = -4108

$ws.Range("J7").Select() | Out-Null
